$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Settings sheet
# ---------------------------------------------------------------------------
$settings = $wb.Worksheets.Item("Settings")
$settings.Range("B2").Value = "P004_SP004_090_NHC_NHP_Performer_Queue"
$settings.Range("B3").Value = "DEV"
$settings.Range("B12").Select()

# ---------------------------------------------------------------------------
# Assets sheet
# ---------------------------------------------------------------------------
$assets = $wb.Worksheets.Item("Assets")
$assets.Range("A2").Value = "AdobeSignEmailAddress"
$assets.Range("B2").Value = "P004_SP004_090_AdobeSignEmailAddress"
$assets.Columns.Item(2).ColumnWidth = 35.77734375
$assets.Range("A2").Select()

# ---------------------------------------------------------------------------
# Constants sheet
# ---------------------------------------------------------------------------
$constants = $wb.Worksheets.Item("Constants")
$constants.Range("A19").Value = "O365AppID"
$constants.Range("B19").Value = "Shared_O365ApplicationID"
$constants.Range("A20").Value = "O365TenantID"
$constants.Range("B20").Value = "Shared_O365TenantID"
$constants.Range("A21").Value = "O365ApplicationSecret"
$constants.Range("B21").Value = "Shared_O365ApplicationSecret"
$constants.Range("A22").Value = "NHC_SharepointURL"
$constants.Range("B22").Value = "https://officemgmtentserv.sharepoint.com/sites/NewHireCommunication"
$constants.Range("A23").Value = "NHC_SharepointListName"
$constants.Range("B23").Value = "New Hire Employee Details"
$constants.Range("A24").Value = "DOH_ColumnDisplayName"
$constants.Range("B24").Value = "DateOfHire"
$constants.Range("A22:B24").Select()
